# entscheid_zustaendigkeit.docx - rsta template changes
#
# 1. "Unsere Referenz" -> "eBau Nummer" (table cell)
# 2. "dem / der Gesuchsteller/in" -> "dem/der Gesuchsteller/in"
# 3. w:overflowPunct flipped from true (default/omitted) to explicit false,
#    everywhere it occurs (document.xml paragraphs + styles.xml style defs)
# 4. A new empty "continuous" section-break paragraph pair inserted near the
#    end of the body (between the existing pair of such paragraphs)
# 5. Nine new character styles ListLabel54..ListLabel62 added to styles.xml,
#    right after the existing ListLabel53 style

$d = $word.ActiveDocument

# ---- simple text fixes (safe via Find/Replace - unique occurrences) ----
$d.Content.Find.Execute("Unsere Referenz", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "eBau Nummer", 2) | Out-Null

$d.Content.Find.Execute("dem / der Gesuchsteller/in", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "dem/der Gesuchsteller/in", 2) | Out-Null

# ---- structural edits: pull the whole package as OOXML, edit the text, push it back ----
$xml = $d.Content.WordOpenXML

# 3. every <w:overflowPunct/> (value defaults to "true") becomes explicit "false"
$xml = $xml.Replace("<w:overflowPunct/>", "<w:overflowPunct w:val=`"false`"/>")

# 4. insert a new (sectPr-continuous paragraph, empty-paragraph) pair right
#    before the last of the two existing "continuous" section-break paragraphs
#    at the tail of the body.
$sectParaPattern = '<w:p w14:paraId="[0-9A-F]+" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:pPr><w:sectPr><w:footnotePr><w:numFmt w:val="decimal"/></w:footnotePr><w:type w:val="continuous"/><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1418" w:right="851" w:bottom="851" w:left="1588" w:header="397" w:footer="0" w:gutter="0"/><w:formProt w:val="0"/><w:textDirection w:val="lrTb"/><w:docGrid w:type="default" w:linePitch="312" w:charSpace="0"/></w:sectPr></w:pPr></w:p>'

$sectMatches = [regex]::Matches($xml, $sectParaPattern)
if ($sectMatches.Count -lt 2) {
    throw "expected at least 2 continuous-sectPr paragraphs, found $($sectMatches.Count)"
}
$lastSectMatch = $sectMatches[$sectMatches.Count - 1]

$newSectPara = '<w:p w14:paraId="0000003B" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:pPr><w:sectPr><w:footnotePr><w:numFmt w:val="decimal"/></w:footnotePr><w:type w:val="continuous"/><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1418" w:right="851" w:bottom="851" w:left="1588" w:header="397" w:footer="0" w:gutter="0"/><w:formProt w:val="0"/><w:textDirection w:val="lrTb"/><w:docGrid w:type="default" w:linePitch="312" w:charSpace="0"/></w:sectPr></w:pPr></w:p>'
$newEmptyPara = '<w:p w14:paraId="0000003C" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Normal"/><w:widowControl/><w:suppressAutoHyphens w:val="0"/><w:overflowPunct w:val="false"/><w:bidi w:val="0"/><w:spacing w:before="0" w:after="220" w:line="280" w:lineRule="atLeast"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>'
$insertion = $newSectPara + $newEmptyPara

$insertAt = $lastSectMatch.Index
$xml = $xml.Substring(0, $insertAt) + $insertion + $xml.Substring($insertAt)

# 5. add ListLabel54..ListLabel62 character styles right after ListLabel53
$listLabel53Anchor = '<w:style w:type="character" w:styleId="ListLabel53"><w:name w:val="ListLabel 53"/><w:qFormat/><w:rPr><w:rFonts w:cs="Wingdings"/></w:rPr></w:style>'
if ($xml.IndexOf($listLabel53Anchor) -lt 0) {
    throw "ListLabel53 style anchor not found"
}

$newLabels = @(
    @{ Num = 54; Font = "Arial" },
    @{ Num = 55; Font = "Courier New" },
    @{ Num = 56; Font = "Wingdings" },
    @{ Num = 57; Font = "Symbol" },
    @{ Num = 58; Font = "Courier New" },
    @{ Num = 59; Font = "Wingdings" },
    @{ Num = 60; Font = "Symbol" },
    @{ Num = 61; Font = "Courier New" },
    @{ Num = 62; Font = "Wingdings" }
)
$labelBlock = ""
foreach ($lbl in $newLabels) {
    $labelBlock += '<w:style w:type="character" w:styleId="ListLabel' + $lbl.Num + '"><w:name w:val="ListLabel ' + $lbl.Num + '"/><w:qFormat/><w:rPr><w:rFonts w:cs="' + $lbl.Font + '"/></w:rPr></w:style>'
}

$xml = $xml.Replace($listLabel53Anchor, $listLabel53Anchor + $labelBlock)

$d.Content.WordOpenXML = $xml

Write-Host "overflowPunct fixed, paragraphs inserted, ListLabel54-62 added, text fixed"
